$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New average of J column
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertical-centered style on a scratch cell, then
# copy that formatting (as a single atomic paste) onto the summary values
# so we don't leave a trail of intermediate cell formats behind.
$scratch = $ws.Range("AB1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Range("A14:B17").Select()
